# FB_012 - update booking report rows 2-6 with new address/user data and
# refreshed totals (proxy module now surfaces rent + product controller
# data, so the sample export picks up new users/addresses/amounts).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Омск Кордная 19 / testUser
$ws.Cells.Item(2, 2).Value = 6
$ws.Cells.Item(2, 3).Value = 15504
$ws.Cells.Item(2, 5).Value = "Омск Кордная 19"
$ws.Cells.Item(2, 6).Value = "testUser"
$ws.Cells.Item(2, 7).Value = 15
$ws.Cells.Item(2, 9).ClearContents()

# Row 3: Краснодар Столовая 3 / testUser
$ws.Cells.Item(3, 2).Value = 372
$ws.Cells.Item(3, 3).Value = 1145016
$ws.Cells.Item(3, 5).Value = "Краснодар Столовая 3"
$ws.Cells.Item(3, 6).Value = "testUser"
$ws.Cells.Item(3, 7).Value = 10
$ws.Cells.Item(3, 9).ClearContents()

# Row 4: Пенза Рабочая 17 / testUser
$ws.Cells.Item(4, 2).Value = 372
$ws.Cells.Item(4, 3).Value = 699732
$ws.Cells.Item(4, 5).Value = "Пенза Рабочая 17"
$ws.Cells.Item(4, 6).Value = "testUser"
$ws.Cells.Item(4, 7).Value = 12
$ws.Cells.Item(4, 9).ClearContents()

# Row 5: Тверь Тверская 15 / testUser
$ws.Cells.Item(5, 2).Value = 372
$ws.Cells.Item(5, 3).Value = 590884.8
$ws.Cells.Item(5, 5).Value = "Тверь Тверская 15"
$ws.Cells.Item(5, 6).Value = "testUser"
$ws.Cells.Item(5, 7).Value = 12
$ws.Cells.Item(5, 9).ClearContents()

# Row 6: Краснодар Столовая 3 / testUser
$ws.Cells.Item(6, 2).Value = 372
$ws.Cells.Item(6, 3).Value = 1119571.2
$ws.Cells.Item(6, 5).Value = "Краснодар Столовая 3"
$ws.Cells.Item(6, 6).Value = "testUser"
$ws.Cells.Item(6, 7).Value = 12
$ws.Cells.Item(6, 9).ClearContents()

Write-Host "applied rows 2-6 updates"
